# Apply commit "feat: add 2022-Q1 data"
#
# Before:
#   Sheet 1 "2021-Q4"  -> fund holdings for 2021-Q4 (unchanged)
#   Sheet 2 "总计"      -> totals summary (date / count / market value)
#
# After:
#   Sheet 1 "2021-Q4"  -> unchanged
#   Sheet 2 "2022-Q1"  -> fund holdings for 2022-Q1 (replaces old "总计" content)
#   Sheet 3 "总计"      -> totals summary, with a new row for 2022-Q1 prepended

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)      # "2021-Q4" - stays untouched
$ws2 = $wb.Worksheets.Item(2)      # currently "总计" -> becomes "2022-Q1"

# --- Add the brand-new "总计" sheet at the end of the workbook ---
$wsTotal = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# match the page-margin setup used by the rest of the workbook (0.75in /
# 0.75in / 1in / 1in / 0.5in / 0.5in, expressed here in points)
$wsTotal.PageSetup.LeftMargin   = 54
$wsTotal.PageSetup.RightMargin  = 54
$wsTotal.PageSetup.TopMargin    = 72
$wsTotal.PageSetup.BottomMargin = 72
$wsTotal.PageSetup.HeaderMargin = 36
$wsTotal.PageSetup.FooterMargin = 36

# Write a value as real text (so things like leading zeros / decimal-looking
# strings such as "007139" or "12.79" are preserved instead of being coerced
# to numbers by Excel) and leave the cell on the default "Normal" style.
function Set-PlainTextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Write a value as real text, then stamp the cell with the same number
# format / font / border / alignment as $formatSource (a single styled
# reference cell), without disturbing the text we just wrote.
function Set-StyledTextValue($range, [string]$text, $formatSource) {
    $range.Value = "'" + $text
    $formatSource.Copy() | Out-Null
    $range.PasteSpecial($xlPasteFormats)
}

# ============================================================
# Sheet 2: rename "总计" -> "2022-Q1" and replace its contents
# ============================================================
$ws2.Name = "2022-Q1"

$headerFmt = $ws1.Cells.Item(1,2)   # styled header cell (bold/border/center)
$indexFmt  = $ws1.Cells.Item(2,1)   # styled index cell in column A

Set-StyledTextValue $ws2.Cells.Item(1,2) "基金代码" $headerFmt
Set-StyledTextValue $ws2.Cells.Item(1,3) "基金名称" $headerFmt
Set-StyledTextValue $ws2.Cells.Item(1,4) "基金规模" $headerFmt
Set-StyledTextValue $ws2.Cells.Item(1,5) "股票总仓位" $headerFmt
Set-StyledTextValue $ws2.Cells.Item(1,6) "仓位占比" $headerFmt
Set-StyledTextValue $ws2.Cells.Item(1,7) "持有市值(亿元)" $headerFmt
Set-StyledTextValue $ws2.Cells.Item(1,8) "仓位排名" $headerFmt

$ws2.Cells.Item(2,1).Value = 0
$indexFmt.Copy() | Out-Null
$ws2.Cells.Item(2,1).PasteSpecial($xlPasteFormats)

Set-PlainTextValue $ws2.Cells.Item(2,2) "007139"
Set-PlainTextValue $ws2.Cells.Item(2,3) "富国民裕进取沪港深成长精选混合"
Set-PlainTextValue $ws2.Cells.Item(2,4) "12.79"
Set-PlainTextValue $ws2.Cells.Item(2,5) "92.21"
Set-PlainTextValue $ws2.Cells.Item(2,6) "5.96"
Set-PlainTextValue $ws2.Cells.Item(2,7) "0.7623"
$ws2.Cells.Item(2,8).Value = 4

# ============================================================
# New sheet: "总计" with the running totals table
# ============================================================
$wsTotal.Name = "总计"

Set-StyledTextValue $wsTotal.Cells.Item(1,2) "日期" $headerFmt
Set-StyledTextValue $wsTotal.Cells.Item(1,3) "持有数量(只)" $headerFmt
Set-StyledTextValue $wsTotal.Cells.Item(1,4) "持有市值(亿元)" $headerFmt

$wsTotal.Cells.Item(2,1).Value = 0
$indexFmt.Copy() | Out-Null
$wsTotal.Cells.Item(2,1).PasteSpecial($xlPasteFormats)
Set-PlainTextValue $wsTotal.Cells.Item(2,2) "2022-Q1"
$wsTotal.Cells.Item(2,3).Value = 1
$wsTotal.Cells.Item(2,4).Value = 0.76

$wsTotal.Cells.Item(3,1).Value = 1
$indexFmt.Copy() | Out-Null
$wsTotal.Cells.Item(3,1).PasteSpecial($xlPasteFormats)
Set-PlainTextValue $wsTotal.Cells.Item(3,2) "2021-Q4"
$wsTotal.Cells.Item(3,3).Value = 1
$wsTotal.Cells.Item(3,4).Value = 0.6
